$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: fill in B18 report text (keeps existing style s="1")
$ws.Range("B18").Value = "The mini-cog was performed on 11/09/2014 and the results were as follows: Recall: 0/3 Clock Drawing Test:  3/4 "

# Row 19: long-integer NoteID in A19 (numeric "0" format) and report text in B19 (no special style)
$ws.Range("A19").Value = 1234567891012
$ws.Range("A19").NumberFormat = "0"
$ws.Range("B19").Value = "..although he did well 26/30 on SLUMS. Loss of weight ~ 2 weeks 175 --> 156 lbs..."
$ws.Range("B19").Style = "Normal"
$ws.Rows.Item(19).EntireRow.AutoFit()

# Row 20: NoteID in A20 (default/general style) and report text in B20 (keeps existing style s="1")
$ws.Range("A20").Value = 1234
$ws.Range("B20").Value = "...VS taken after SLUMS and GDS BS R 162/80, L 168/84…"

# Update the active selection to D17
[void]$ws.Range("D17").Select()
